$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.6506806930073017, 0.9895372330898513, 0.3895017181495386, 2.092490195104551)"
$ws.Range("C2").Value = "NIG(0.6416845190653832, 0.441633232529803, 5.644212749393322, 4.71719442309704)"
$ws.Range("D2").Value = "F(305.4486421138955, 16.121714900168392, -9.136066842126588, 12.328917391663836)"
$ws.Range("E2").Value = "NIG(1.4737684713430976, 1.1594540624311656, 4.179205171568178, 6.110520539946036)"
